# Refresh the cryptos price/volume snapshot (and fix the TRON/ShibaInu
# row order swap) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '27.368.67'
$ws.Range("E2").Value2 = '  -2.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '1.857.94'
$ws.Range("E3").Value2 = '  -3.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '1.002'
$ws.Range("E4").Value2 = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '328.22'
$ws.Range("E5").Value2 = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '1.002'
$ws.Range("E6").Value2 = '  +0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.4601'
$ws.Range("E7").Value2 = '  -1.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.3935'
$ws.Range("E8").Value2 = '  -1.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '46.91'
$ws.Range("E9").Value2 = '  -11.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '0.07930'
$ws.Range("E10").Value2 = '  -5.54%  '
$ws.Range("E11").Value2 = '  -3.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '21.51'
$ws.Range("E12").Value2 = '  -2.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '1.857.21'
$ws.Range("E13").Value2 = '  -1.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '5.931'
$ws.Range("E14").Value2 = '  -2.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '7.153'
$ws.Range("E15").Value2 = '  -3.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '1.001'
$ws.Range("E16").Value2 = '  +0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '86.47'
$ws.Range("E17").Value2 = '  -3.62%  '
$ws.Range("B18").Value2 = 'ShibaInu'
$ws.Range("C18").Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '0.00001031'
$ws.Range("E18").Value2 = '  -2.73%  '
$ws.Range("B19").Value2 = 'TRON'
$ws.Range("C19").Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '0.06575'
$ws.Range("E19").Value2 = '  -0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '17.22'
$ws.Range("E20").Value2 = '  -4.13%  '
$ws.Range("E21").Value2 = '  +0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '5.493'
$ws.Range("E22").Value2 = '  -4.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '27.387.31'
$ws.Range("E23").Value2 = '  -2.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '10.92'
$ws.Range("E24").Value2 = '  -3.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '2.305'
$ws.Range("E25").Value2 = '  +1.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '2.078.63'
$ws.Range("E26").Value2 = '  -1.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '153.65'
$ws.Range("E27").Value2 = '  +0.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '20.09'
$ws.Range("E28").Value2 = '  +0.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '2.071'
$ws.Range("E29").Value2 = '  -2.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '5.483'
$ws.Range("E30").Value2 = '  -4.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '121.59'
$ws.Range("E31").Value2 = '  -1.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '0.9555'
$ws.Range("E32").Value2 = '  -2.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '0.09411'
$ws.Range("E33").Value2 = '  -2.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '1.453'
$ws.Range("E34").Value2 = '  +0.35%  '
$ws.Range("E35").Value2 = '  -1.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '5.278'
$ws.Range("E36").Value2 = '  -4.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '0.06041'
$ws.Range("E37").Value2 = '  -2.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '0.02231'
$ws.Range("E38").Value2 = '  -3.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '1.221'
$ws.Range("E39").Value2 = '  -2.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '8.080'
$ws.Range("E40").Value2 = '  -8.25%  '
$ws.Range("E41").Value2 = '  +0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '0.5935'
$ws.Range("E42").Value2 = '  -3.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.1895'
$ws.Range("E43").Value2 = '  -0.46%  '
$ws.Range("E44").Value2 = '  -7.83%  '
$ws.Range("E45").Value2 = '  -1.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '0.5642'
$ws.Range("E46").Value2 = '  -3.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '12.05'
$ws.Range("E47").Value2 = '  -5.28%  '
$ws.Range("E48").Value2 = '  -1.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '1.921'
$ws.Range("E49").Value2 = '  -5.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '0.06756'
$ws.Range("E50").Value2 = '  -1.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '108.72'
$ws.Range("E51").Value2 = '  -1.13%  '
